$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2039106145251397
$ws.Range("C2").Value = 0.5223463687150838
$ws.Range("J2").Value = 0.02793296089385475
$ws.Range("P2").Value = 0.1592178770949721
$ws.Range("S2").Value = 0.08659217877094973
$ws.Range("B3").Value = 0.005154639175257732
$ws.Range("C3").Value = 0.02577319587628866
$ws.Range("J3").Value = 0.05670103092783505
$ws.Range("P3").Value = 0.7268041237113402
$ws.Range("S3").Value = 0.1855670103092784
$ws.Range("J4").Value = 0.06976744186046512
$ws.Range("P4").Value = 0.5813953488372093
$ws.Range("S4").Value = 0.3488372093023256
$ws.Range("B6").Value = 0.1238095238095238
$ws.Range("D6").Value = 0.009523809523809525
$ws.Range("F6").Value = 0.04761904761904762
$ws.Range("J6").Value = 0.2476190476190476
$ws.Range("O6").Value = 0.03333333333333333
$ws.Range("Q6").Value = 0.1285714285714286
$ws.Range("R6").Value = 0.05238095238095238
$ws.Range("S6").Value = 0.3571428571428572
$ws.Range("B7").Value = 0.1357466063348416
$ws.Range("D7").Value = 0.03167420814479638
$ws.Range("E7").Value = 0.004524886877828055
$ws.Range("F7").Value = 0.05429864253393665
$ws.Range("J7").Value = 0.16289592760181
$ws.Range("O7").Value = 0.02262443438914027
$ws.Range("Q7").Value = 0.1447963800904978
$ws.Range("R7").Value = 0.06787330316742081
$ws.Range("S7").Value = 0.3755656108597285
$ws.Range("B8").Value = 0.08355091383812011
$ws.Range("D8").Value = 0.01305483028720627
$ws.Range("F8").Value = 0.06266318537859007
$ws.Range("J8").Value = 0.1201044386422977
$ws.Range("O8").Value = 0.02088772845953003
$ws.Range("Q8").Value = 0.1462140992167102
$ws.Range("R8").Value = 0.08093994778067885
$ws.Range("S8").Value = 0.4725848563968668
$ws.Range("B9").Value = 0.07253886010362694
$ws.Range("D9").Value = 0.0310880829015544
$ws.Range("F9").Value = 0.03626943005181347
$ws.Range("J9").Value = 0.1968911917098446
$ws.Range("Q9").Value = 0.1761658031088083
$ws.Range("R9").Value = 0.05181347150259067
$ws.Range("S9").Value = 0.4352331606217616
$ws.Range("B10").Value = 0.1340057636887608
$ws.Range("D10").Value = 0.01657060518731988
$ws.Range("E10").Value = 0.001440922190201729
$ws.Range("F10").Value = 0.06772334293948126
$ws.Range("J10").Value = 0.1195965417867435
$ws.Range("O10").Value = 0.01945244956772334
$ws.Range("Q10").Value = 0.2046109510086455
$ws.Range("R10").Value = 0.07564841498559077
$ws.Range("S10").Value = 0.3609510086455331
$ws.Range("G11").Value = 0.1483679525222552
$ws.Range("J11").Value = 0.1008902077151335
$ws.Range("K11").Value = 0.198813056379822
$ws.Range("L11").Value = 0.5311572700296736
$ws.Range("S11").Value = 0.02077151335311573
$ws.Range("G12").Value = 0.7438423645320197
$ws.Range("J12").Value = 0.167487684729064
$ws.Range("K12").Value = 0.02463054187192118
$ws.Range("L12").Value = 0.03940886699507389
$ws.Range("S12").Value = 0.02463054187192118
$ws.Range("G13").Value = 0.8235294117647058
$ws.Range("J13").Value = 0.1176470588235294
$ws.Range("S13").Value = 0.05882352941176471
$ws.Range("F15").Value = 0.04310344827586207
$ws.Range("H15").Value = 0.146551724137931
$ws.Range("I15").Value = 0.06896551724137931
$ws.Range("J15").Value = 0.396551724137931
$ws.Range("K15").Value = 0.03448275862068965
$ws.Range("O15").Value = 0.04741379310344827
$ws.Range("S15").Value = 0.2629310344827586
$ws.Range("F16").Value = 0.01913875598086124
$ws.Range("H16").Value = 0.1818181818181818
$ws.Range("I16").Value = 0.07655502392344497
$ws.Range("J16").Value = 0.4688995215311005
$ws.Range("K16").Value = 0.1052631578947368
$ws.Range("O16").Value = 0.02392344497607655
$ws.Range("S16").Value = 0.1244019138755981
$ws.Range("F17").Value = 0.007009345794392523
$ws.Range("H17").Value = 0.1144859813084112
$ws.Range("I17").Value = 0.1074766355140187
$ws.Range("J17").Value = 0.4742990654205608
$ws.Range("K17").Value = 0.09813084112149532
$ws.Range("M17").Value = 0.01635514018691589
$ws.Range("O17").Value = 0.0514018691588785
$ws.Range("S17").Value = 0.1308411214953271
$ws.Range("F18").Value = 0.005882352941176471
$ws.Range("H18").Value = 0.1352941176470588
$ws.Range("I18").Value = 0.06470588235294118
$ws.Range("J18").Value = 0.4529411764705882
$ws.Range("K18").Value = 0.1176470588235294
$ws.Range("M18").Value = 0.01764705882352941
$ws.Range("N18").Value = 0.005882352941176471
$ws.Range("O18").Value = 0.08823529411764706
$ws.Range("S18").Value = 0.1117647058823529
$ws.Range("F19").Value = 0.01227935533384497
$ws.Range("H19").Value = 0.182655410590944
$ws.Range("I19").Value = 0.07444359171143515
$ws.Range("J19").Value = 0.3898695318495779
$ws.Range("K19").Value = 0.1266308518802763
$ws.Range("M19").Value = 0.01765157329240215
$ws.Range("N19").Value = 0.003837298541826554
$ws.Range("O19").Value = 0.08288564850345356
$ws.Range("S19").Value = 0.1097467382962394
